$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "317.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.71%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.45%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.183"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.10%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08032"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.98%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.491"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.81%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.517"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.92%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.935"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.26%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.000"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.99%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9388"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.23%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1286"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "9.50%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1938"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.21%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09007"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.91%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03393"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.22%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09546"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.57%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001391"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.57%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006179"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.70%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.378"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.43%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3524"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.38%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.592"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "25.15%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1316"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.46%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2427"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.18%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04385"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.59%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.43%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004270"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-8.54%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001329"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.18%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003989"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.07%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02357"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.17%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05153"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.31%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007624"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.70%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1401"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.50%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008688"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.88%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002109"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.21%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008843"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.70%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006485"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.16%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002859"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-6.39%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001689"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "69.06%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.03%"
Write-Output "Updated symbol list values."
